$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "634×7=" "306×9="
Replace-Text "251×6=" "699×2="
Replace-Text "408×4=" "499×3="
Replace-Text "524×4=" "655×4="
Replace-Text "886×2=" "879×6="
Replace-Text "763×5=" "884×4="
Replace-Text "973×4=" "135×9="
Replace-Text "529×6=" "118×3="
Replace-Text "494×7=" "935×4="
Replace-Text "135×8=" "656×5="
Replace-Text "711×7=" "166×8="
Replace-Text "237×5=" "545×4="
Replace-Text "730×4=" "587×3="
Replace-Text "508×9=" "346×5="
Replace-Text "117×8=" "750×6="
Replace-Text "748×3=" "813×7="
Replace-Text "722×5=" "859×8="
Replace-Text "243×8=" "334×9="
Replace-Text "509×9=" "907×2="
Replace-Text "902×5=" "492×9="
Replace-Text "870×4=" "860×5="
Replace-Text "715×8=" "951×5="
Replace-Text "932×6=" "206×7="
Replace-Text "921×2=" "965×6="
Replace-Text "357×9=" "886×9="

Write-Output "Done"
